# Regenerated staging-template header row: a new "BusinessKey" column is
# inserted as the first data column on row 2 (ahead of the existing
# Code / DataVersion_ID / Description / Name / Order columns), which
# shifts those five existing header labels one column to the right and
# adds "Order" in the new, sixth (F) column.
#
#   before: A2=Code  B2=DataVersion_ID  C2=Description  D2=Name  E2=Order
#   after:  A2=BusinessKey  B2=Code  C2=DataVersion_ID  D2=Description  E2=Name  F2=Order

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new F2 cell the same (bold + underlined) header formatting as
# the rest of row 2 before populating it, using the existing E2 header
# cell as the format source.
$ws.Range("E2").Copy()
$ws.Range("F2").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Shift the existing header labels one column to the right, working from
# the right-most column back towards A so no value is clobbered before
# it has been copied onward.
$ws.Range("F2").Value2 = $ws.Range("E2").Value2  # Order
$ws.Range("E2").Value2 = $ws.Range("D2").Value2  # Name
$ws.Range("D2").Value2 = $ws.Range("C2").Value2  # Description
$ws.Range("C2").Value2 = $ws.Range("B2").Value2  # DataVersion_ID
$ws.Range("B2").Value2 = $ws.Range("A2").Value2  # Code

# Populate the freed-up first column with the new field name.
$ws.Range("A2").Value2 = "BusinessKey"
